# Usability Testing 4 & 5 notes and updates
# Adds "Person 5" (col G) and "Person 6" (col H) to the Notes&Themes sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Notes&Themes")
$ws2 = $wb.Worksheets.Item("Feature Importance")

# --- Copy cell formatting (fill/border/alignment) from column F into the
# --- new columns G (Person 5) and H (Person 6) before filling in values,
# --- so the new cells match the existing header/body styling exactly.
$ws1.Range("F2").Copy() | Out-Null
$ws1.Range("G2:H2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws1.Range("F3:F12").Copy() | Out-Null
$ws1.Range("G3:G12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws1.Range("F3:F12").Copy() | Out-Null
$ws1.Range("H3:H12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Header row ---
$ws1.Range("G2").Value = "Person 5"
$ws1.Range("H2").Value = "Person 6"

# --- Person 5 (Usability Testing 4) ---
$ws1.Range("G3").Value = "12 & 9"
$ws1.Range("G4").Value = "Yes, one child has a heart condition"
$ws1.Range("G5").Value = "usual, Allergies, headaches"
$ws1.Range("G6").Value = """just figuring it out"""
$ws1.Range("G7").Value = "Which drugs hasve interactions with her child's heart medication."
$ws1.Range("G8").Value = "History, saving notes about past experiences with medicinesm (e.g. which allergy medicine has worked best)"
$ws1.Range("G9").Value = "Like an actual cabinet, user profiles/saving medicines and searches"
$ws1.Range("G10").Value = "would like to search by active ingredient"
$ws1.Range("G11").Value = "would use it mostly to make sure medicines she gives her son and family would not hurt them or interact"
$ws1.Range("G12").Value = "would like to write notes and keep a history"

# --- Person 6 (Usability Testing 5) ---
$ws1.Range("H3").Value = "3, 1.5"
$ws1.Range("H4").Value = "NA"
$ws1.Range("H5").Value = "NA"
$ws1.Range("H6").Value = "Wife is a NICU nurse, she handles all of the meds and healthcare"
$ws1.Range("H7").Value = "A backup or secondary check"
$ws1.Range("H8").Value = "Would like to see what medicines and at what time family members need to take medicine as a reminder"
$ws1.Range("H9").Value = "notifications/reminder"
$ws1.Range("H10").Value = "NA"
$ws1.Range("H11").Value = "Would never use the app since wife is a healthcare professional"
$ws1.Range("H12").Value = "NA"

# --- Column widths for the two new columns (match existing 37-char columns) ---
$ws1.Range("G1:H1").EntireColumn.ColumnWidth = 36.17

# --- Group/outline the original four interviewee columns (C:F) like the
# --- older columns are now collapsible behind the two newest interviews ---
$ws1.Range("C1:F1").EntireColumn.OutlineLevel = 1

# --- Row heights grew to fit the additional wrapped text ---
$ws1.Rows.Item(6).RowHeight = 75
$ws1.Rows.Item(7).RowHeight = 30
$ws1.Rows.Item(8).RowHeight = 60
$ws1.Rows.Item(9).RowHeight = 30
$ws1.Rows.Item(10).RowHeight = 30
$ws1.Rows.Item(11).RowHeight = 60
$ws1.Rows.Item(12).RowHeight = 75

# --- Selection / active sheet bookkeeping ---
$ws2.Range("G11").Select()
$ws1.Activate()
$ws1.Range("J16").Select()
$excel.ActiveWindow.Zoom = 80
